$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 876.7708
$ws.Cells.Item(6, 9).Value = 108.382355
$ws.Cells.Item(6, 10).Value = 2742.8572
$ws.Cells.Item(6, 11).Value = 325.147065
$ws.Cells.Item(6, 12).Value = 8228.571599999999
$ws.Cells.Item(6, 13).Value = -213.147065
$ws.Cells.Item(6, 14).Value = -8452.571599999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = $null
$ws.Cells.Item(12, 14).Value = $null

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 16883.777
$ws.Cells.Item(28, 9).Value = 278.2857
$ws.Cells.Item(28, 11).Value = 278.2857
$ws.Cells.Item(28, 13).Value = 206.7143

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 702.7568
$ws.Cells.Item(33, 9).Value = 529.0357
$ws.Cells.Item(33, 10).Value = 1243.2222
$ws.Cells.Item(33, 11).Value = 529.0357
$ws.Cells.Item(33, 12).Value = 1243.2222
$ws.Cells.Item(33, 13).Value = -300.0357
$ws.Cells.Item(33, 14).Value = -1701.2222

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2821
$ws.Cells.Item(62, 9).Value = 2152.5
$ws.Cells.Item(62, 10).Value = 3266.6667
$ws.Cells.Item(62, 11).Value = 2152.5
$ws.Cells.Item(62, 12).Value = 3266.6667
$ws.Cells.Item(62, 13).Value = -1528.5
$ws.Cells.Item(62, 14).Value = -4514.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 2821
$ws.Cells.Item(65, 9).Value = 2152.5
$ws.Cells.Item(65, 10).Value = 3266.6667
$ws.Cells.Item(65, 11).Value = 10762.5
$ws.Cells.Item(65, 12).Value = 16333.3335
$ws.Cells.Item(65, 13).Value = -7642.5
$ws.Cells.Item(65, 14).Value = -22573.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 3739.1316
$ws.Cells.Item(132, 9).Value = 3735.3713
$ws.Cells.Item(132, 10).Value = 3783
$ws.Cells.Item(132, 11).Value = 11206.1139
$ws.Cells.Item(132, 12).Value = 11349
$ws.Cells.Item(132, 13).Value = -8676.1139
$ws.Cells.Item(132, 14).Value = -16409

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 3090.6428
$ws.Cells.Item(137, 9).Value = 2752.6667
$ws.Cells.Item(137, 10).Value = 4329.8887
$ws.Cells.Item(137, 11).Value = 8258.000100000001
$ws.Cells.Item(137, 12).Value = 12989.6661
$ws.Cells.Item(137, 13).Value = -5708.000100000001
$ws.Cells.Item(137, 14).Value = -18089.6661

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2035.5834
$ws.Cells.Item(138, 9).Value = 2246.1
$ws.Cells.Item(138, 11).Value = 6738.299999999999
$ws.Cells.Item(138, 13).Value = -1598.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1905.25
$ws.Cells.Item(2, 9).Value = 1905.25
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 1905.25
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -1792.25
$ws.Cells.Item(2, 14).Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 555367.8
$ws.Cells.Item(32, 9).Value = 651259.4399999999
$ws.Cells.Item(32, 10).Value = 22636.445
$ws.Cells.Item(32, 11).Value = 651259.4399999999
$ws.Cells.Item(32, 12).Value = 22636.445
$ws.Cells.Item(32, 13).Value = -650972.4399999999
$ws.Cells.Item(32, 14).Value = -23210.445

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(40, 8).Value = 5959.75
$ws.Cells.Item(40, 10).Value = 5959.75
$ws.Cells.Item(40, 12).Value = 5959.75
$ws.Cells.Item(40, 14).Value = -6311.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2475.6667
$ws.Cells.Item(45, 9).Value = 1787.1578
$ws.Cells.Item(45, 11).Value = 1787.1578
$ws.Cells.Item(45, 13).Value = -1410.1578

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3694.923
$ws.Cells.Item(61, 9).Value = 3319.1428
$ws.Cells.Item(61, 10).Value = 4133.3335
$ws.Cells.Item(61, 11).Value = 3319.1428
$ws.Cells.Item(61, 12).Value = 4133.3335
$ws.Cells.Item(61, 13).Value = -3107.1428
$ws.Cells.Item(61, 14).Value = -4557.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 1905.25
$ws.Cells.Item(116, 9).Value = 1905.25
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 1905.25
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = 388.75
$ws.Cells.Item(116, 14).Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 6945.4
$ws.Cells.Item(132, 9).Value = 9402.4
$ws.Cells.Item(132, 10).Value = 4488.4
$ws.Cells.Item(132, 11).Value = 28207.2
$ws.Cells.Item(132, 12).Value = 13465.2
$ws.Cells.Item(132, 13).Value = -25677.2
$ws.Cells.Item(132, 14).Value = -18525.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 3694.923
$ws.Cells.Item(136, 9).Value = 3319.1428
$ws.Cells.Item(136, 10).Value = 4133.3335
$ws.Cells.Item(136, 11).Value = 9957.428400000001
$ws.Cells.Item(136, 12).Value = 12400.0005
$ws.Cells.Item(136, 13).Value = -7407.428400000001
$ws.Cells.Item(136, 14).Value = -17500.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1905.25
$ws.Cells.Item(3, 9).Value = 1905.25
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 1905.25
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = -1791.25
$ws.Cells.Item(3, 14).Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 12502397
$ws.Cells.Item(105, 9).Value = 13891174
$ws.Cells.Item(105, 10).Value = 3400
$ws.Cells.Item(105, 11).Value = 13891174
$ws.Cells.Item(105, 12).Value = 3400
$ws.Cells.Item(105, 13).Value = -13889427
$ws.Cells.Item(105, 14).Value = -6894

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(109, 8).Value = 90000
$ws.Cells.Item(109, 10).Value = 90000
$ws.Cells.Item(109, 12).Value = 90000
$ws.Cells.Item(109, 14).Value = -92774

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7015.3335
$ws.Cells.Item(31, 9).Value = 1272.3636
$ws.Cells.Item(31, 10).Value = 18501.273
$ws.Cells.Item(31, 11).Value = 1272.3636
$ws.Cells.Item(31, 12).Value = 18501.273
$ws.Cells.Item(31, 13).Value = -977.3635999999999
$ws.Cells.Item(31, 14).Value = -19091.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 7015.3335
$ws.Cells.Item(34, 9).Value = 1272.3636
$ws.Cells.Item(34, 10).Value = 18501.273
$ws.Cells.Item(34, 11).Value = 1272.3636
$ws.Cells.Item(34, 12).Value = 18501.273
$ws.Cells.Item(34, 13).Value = -1070.3636
$ws.Cells.Item(34, 14).Value = -18905.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 17165.834
$ws.Cells.Item(41, 9).Value = 5000
$ws.Cells.Item(41, 10).Value = 19599
$ws.Cells.Item(41, 11).Value = 5000
$ws.Cells.Item(41, 12).Value = 19599
$ws.Cells.Item(41, 13).Value = -4572
$ws.Cells.Item(41, 14).Value = -20455

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(53, 8).Value = 28450
$ws.Cells.Item(53, 10).Value = 28450
$ws.Cells.Item(53, 12).Value = 28450
$ws.Cells.Item(53, 14).Value = -29664

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 932.625
$ws.Cells.Item(105, 9).Value = 1012.5
$ws.Cells.Item(105, 10).Value = 852.75
$ws.Cells.Item(105, 11).Value = 1012.5
$ws.Cells.Item(105, 12).Value = 852.75
$ws.Cells.Item(105, 13).Value = 734.5
$ws.Cells.Item(105, 14).Value = -4346.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1014.44446
$ws.Cells.Item(5, 9).Value = 890
$ws.Cells.Item(5, 10).Value = 1450
$ws.Cells.Item(5, 11).Value = 2670
$ws.Cells.Item(5, 12).Value = 4350
$ws.Cells.Item(5, 13).Value = -2558
$ws.Cells.Item(5, 14).Value = -4574

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(76, 8).Value = 2785.7144
$ws.Cells.Item(76, 9).Value = 750
$ws.Cells.Item(76, 10).Value = 3000
$ws.Cells.Item(76, 11).Value = 2250
$ws.Cells.Item(76, 12).Value = 9000
$ws.Cells.Item(76, 13).Value = -1867
$ws.Cells.Item(76, 14).Value = -9766

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(79, 8).Value = 2785.7144
$ws.Cells.Item(79, 9).Value = 750
$ws.Cells.Item(79, 10).Value = 3000
$ws.Cells.Item(79, 11).Value = 2250
$ws.Cells.Item(79, 12).Value = 9000
$ws.Cells.Item(79, 13).Value = -924
$ws.Cells.Item(79, 14).Value = -11652

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 758.6053000000001
$ws.Cells.Item(113, 9).Value = 449.125
$ws.Cells.Item(113, 10).Value = 1289.1428
$ws.Cells.Item(113, 11).Value = 1347.375
$ws.Cells.Item(113, 12).Value = 3867.4284
$ws.Cells.Item(113, 13).Value = 822.625
$ws.Cells.Item(113, 14).Value = -8207.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 17056.5
$ws.Cells.Item(122, 9).Value = 468
$ws.Cells.Item(122, 11).Value = 4212
$ws.Cells.Item(122, 13).Value = -1762

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 1014.44446
$ws.Cells.Item(135, 9).Value = 890
$ws.Cells.Item(135, 10).Value = 1450
$ws.Cells.Item(135, 11).Value = 8010
$ws.Cells.Item(135, 12).Value = 13050
$ws.Cells.Item(135, 13).Value = -5475
$ws.Cells.Item(135, 14).Value = -18120

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 80006.75
$ws.Cells.Item(22, 10).Value = 80006.75
$ws.Cells.Item(22, 12).Value = 80006.75
$ws.Cells.Item(22, 14).Value = -81064.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 14).Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 4149.951
$ws.Cells.Item(122, 9).Value = 1433.8667
$ws.Cells.Item(122, 10).Value = 5716.923
$ws.Cells.Item(122, 11).Value = 4301.6001
$ws.Cells.Item(122, 12).Value = 17150.769
$ws.Cells.Item(122, 13).Value = -1851.6001
$ws.Cells.Item(122, 14).Value = -22050.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(42, 8).Value = 70049
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 13).Value = $null
